# Updated symbol list on Wed Feb 15 23:36:44 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) / Volume(1h) (E) columns for this run's crypto snapshot.
# The sheet stores every data value as literal text (the source writer never used
# numeric cells), so each write explicitly formats the cell as Text first -
# otherwise a numeric-looking string like "316.07" or "6.62%" would be auto-
# converted into a real number/percentage by Excel - and resets the style back to
# Normal afterwards so the cell keeps the workbook's original unstyled look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
$newValues = @("316.07", "6.62%", "45.31", "7.41%", "5.174", "2.94%", "0.08086", "6.52%", "4.536", "1.676", "4.00%", "17.30%", "0.1305", "7.89%", "0.1934", "5.08%", "0.09450", "5.15%", "0.04243", "5.72%", "0.1044", "-0.62%", "0.001344", "5.06%", "0.005880", "1.26%", "1.01%", "0.23%", "0.3371", "1.53%", "8.227", "4.33%", "0.1386", "0.04263", "4.94%", "0.001281", "1.16%", "0.004217", "7.55%", "9.37%", "0.02705", "11.77%", "0.05465", "4.83%", "0.005864", "-3.06%", "0.007756", "-0.57%", "6.62%", "0.007376", "-2.11%", "0.008587", "18.49%", "5.51%", "0.00006801", "0.36%", "0.00000000748", "-0.32%", "53.78%", "0.003988", "-5.07%", "0.00002094", "-0.32%", "0.0001994", "-0.32%")

for ($i = 0; $i -lt $cells.Count; $i++) {
    $rng = $ws.Range($cells[$i])
    $rng.NumberFormat = "@"
    $rng.Value = $newValues[$i]
    $rng.Style = "Normal"
}

Write-Output "Updated $($cells.Count) cells"
